$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 173.75  # was 192.14285
$ws.Range("I12").Value = 198.5  # was 229.2
$ws.Range("K12").Value = 198.5  # was 229.2
$ws.Range("M12").Value = -28.5  # was -59.19999999999999

$ws.Range("H64").Value = 2500  # was 0
$ws.Range("J64").Value = 2500  # was 0
$ws.Range("L64").Value = 2500  # was 0
$ws.Range("N64").Value = -2996  # was None

$ws.Range("H67").Value = 2500  # was 0
$ws.Range("J67").Value = 2500  # was 0
$ws.Range("L67").Value = 2500  # was 0
$ws.Range("N67").Value = -4216  # was None

$ws.Range("H74").Value = 127997  # was 3998.5
$ws.Range("I74").Value = 3996  # was 3999
$ws.Range("J74").Value = 500000  # was 3998
$ws.Range("K74").Value = 3996  # was 3999
$ws.Range("L74").Value = 500000  # was 3998
$ws.Range("M74").Value = -3060  # was -3063
$ws.Range("N74").Value = -501872  # was -5870

$ws.Range("H77").Value = 127997  # was 3998.5
$ws.Range("I77").Value = 3996  # was 3999
$ws.Range("J77").Value = 500000  # was 3998
$ws.Range("K77").Value = 19980  # was 19995
$ws.Range("L77").Value = 2500000  # was 19990
$ws.Range("M77").Value = -15300  # was -15315
$ws.Range("N77").Value = -2509360  # was -29350

$ws.Range("H98").Value = 3241.0557  # was 3320.7058
$ws.Range("I98").Value = 3063.6667  # was 3170.6365
$ws.Range("K98").Value = 3063.6667  # was 3170.6365
$ws.Range("M98").Value = -1565.6667  # was -1672.6365

$ws.Range("H112").Value = 2246.9412  # was 2264.5881
$ws.Range("J112").Value = 2320  # was 2299.9375
$ws.Range("L112").Value = 6960  # was 6899.8125
$ws.Range("N112").Value = -9176  # was -9115.8125

$ws.Range("H122").Value = 3241.0557  # was 3320.7058
$ws.Range("I122").Value = 3063.6667  # was 3170.6365
$ws.Range("K122").Value = 9191.000100000001  # was 9511.9095
$ws.Range("M122").Value = -6741.000100000001  # was -7061.9095

$ws.Range("H132").Value = 1191.2  # was 1189
$ws.Range("I132").Value = 715.1  # was 735.0476
$ws.Range("K132").Value = 2145.3  # was 2205.1428
$ws.Range("M132").Value = 384.6999999999998  # was 324.8571999999999

$ws.Range("H138").Value = 1647.6538  # was 1858.3182
$ws.Range("I138").Value = 473.33334  # was 442
$ws.Range("K138").Value = 1420.00002  # was 1326
$ws.Range("M138").Value = 3719.99998  # was 3814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2494.2727  # was 2638.3333
$ws.Range("I2").Value = 2604.6667  # was 2821.1428
$ws.Range("J2").Value = 1997.5  # was 1998.5
$ws.Range("K2").Value = 2604.6667  # was 2821.1428
$ws.Range("L2").Value = 1997.5  # was 1998.5
$ws.Range("M2").Value = -2491.6667  # was -2708.1428
$ws.Range("N2").Value = -2223.5  # was -2224.5

$ws.Range("H63").Value = 3011.375  # was 2873.25
$ws.Range("I63").Value = 1796  # was 1795.3334
$ws.Range("J63").Value = 3416.5  # was 3520
$ws.Range("K63").Value = 1796  # was 1795.3334
$ws.Range("L63").Value = 3416.5  # was 3520
$ws.Range("M63").Value = -1110  # was -1109.3334
$ws.Range("N63").Value = -4788.5  # was -4892

$ws.Range("H66").Value = 3011.375  # was 2873.25
$ws.Range("I66").Value = 1796  # was 1795.3334
$ws.Range("J66").Value = 3416.5  # was 3520
$ws.Range("K66").Value = 8980  # was 8976.666999999999
$ws.Range("L66").Value = 17082.5  # was 17600
$ws.Range("M66").Value = -5548  # was -5544.666999999999
$ws.Range("N66").Value = -23946.5  # was -24464

$ws.Range("H116").Value = 2494.2727  # was 2638.3333
$ws.Range("I116").Value = 2604.6667  # was 2821.1428
$ws.Range("J116").Value = 1997.5  # was 1998.5
$ws.Range("K116").Value = 2604.6667  # was 2821.1428
$ws.Range("L116").Value = 1997.5  # was 1998.5
$ws.Range("M116").Value = -310.6667000000002  # was -527.1428000000001
$ws.Range("N116").Value = -6585.5  # was -6586.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2494.2727  # was 2638.3333
$ws.Range("I3").Value = 2604.6667  # was 2821.1428
$ws.Range("J3").Value = 1997.5  # was 1998.5
$ws.Range("K3").Value = 2604.6667  # was 2821.1428
$ws.Range("L3").Value = 1997.5  # was 1998.5
$ws.Range("M3").Value = -2490.6667  # was -2707.1428
$ws.Range("N3").Value = -2225.5  # was -2226.5

$ws.Range("H20").Value = 1491.6666  # was 1528.5
$ws.Range("I20").Value = 705  # was 737.7778
$ws.Range("J20").Value = 3065  # was 2951.8
$ws.Range("K20").Value = 705  # was 737.7778
$ws.Range("L20").Value = 3065  # was 2951.8
$ws.Range("M20").Value = -458  # was -490.7778
$ws.Range("N20").Value = -3559  # was -3445.8

$ws.Range("H97").Value = 14464.111  # was 14466.333
$ws.Range("I97").Value = 4413  # was 4419.6665
$ws.Range("K97").Value = 4413  # was 4419.6665
$ws.Range("M97").Value = -3422  # was -3428.6665

$ws.Range("H132").Value = 96999.5  # was 97000
$ws.Range("J132").Value = 96999.5  # was 97000
$ws.Range("L132").Value = 96999.5  # was 97000
$ws.Range("N132").Value = -107119.5  # was -107120

$ws.Range("H134").Value = 5782.4  # was 6008.706
$ws.Range("I134").Value = 5415.5  # was 5743.2
$ws.Range("J134").Value = 7250  # was 8000
$ws.Range("K134").Value = 16246.5  # was 17229.6
$ws.Range("L134").Value = 21750  # was 24000
$ws.Range("M134").Value = -13711.5  # was -14694.6
$ws.Range("N134").Value = -26820  # was -29070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 125  # was 120
$ws.Range("J52").Value = 125  # was 120
$ws.Range("L52").Value = 375  # was 360
$ws.Range("N52").Value = -907  # was -892

$ws.Range("H68").Value = 982.1667  # was 952.7692
$ws.Range("J68").Value = 798.2222  # was 778.4
$ws.Range("L68").Value = 2394.6666  # was 2335.2
$ws.Range("N68").Value = -4016.6666  # was -3957.2

$ws.Range("H71").Value = 982.1667  # was 952.7692
$ws.Range("J71").Value = 798.2222  # was 778.4
$ws.Range("L71").Value = 7183.999800000001  # was 7005.599999999999
$ws.Range("N71").Value = -15295.9998  # was -15117.6

$ws.Range("H81").Value = 11720.571  # was 8674
$ws.Range("J81").Value = 13340.667  # was 10008.8
$ws.Range("L81").Value = 40022.001  # was 30026.4
$ws.Range("N81").Value = -42268.001  # was -32272.4

$ws.Range("H84").Value = 11720.571  # was 8674
$ws.Range("J84").Value = 13340.667  # was 10008.8
$ws.Range("L84").Value = 120066.003  # was 90079.2
$ws.Range("N84").Value = -131298.003  # was -101311.2

$ws.Range("H109").Value = 2053.25  # was 1924.875
$ws.Range("I109").Value = 2308.6667  # was 1966.3334
$ws.Range("K109").Value = 6926.000100000001  # was 5899.0002
$ws.Range("M109").Value = -5886.000100000001  # was -4859.0002

$ws.Range("H132").Value = 2833.3333  # was 1886.1428
$ws.Range("I132").Value = 0  # was 1175.75
$ws.Range("K132").Value = 0  # was 10581.75
$ws.Range("M132").ClearContents()  # was -8051.75

$ws.Range("H137").Value = 4424.75  # was 3875.6667
$ws.Range("I137").Value = 3999  # was 3110.5
$ws.Range("J137").Value = 4566.6665  # was 4258.25
$ws.Range("K137").Value = 11997  # was 9331.5
$ws.Range("L137").Value = 13699.9995  # was 12774.75
$ws.Range("M137").Value = -6897  # was -4231.5
$ws.Range("N137").Value = -23899.9995  # was -22974.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 49999  # was 0
$ws.Range("J64").Value = 49999  # was 0
$ws.Range("L64").Value = 49999  # was 0
$ws.Range("N64").Value = -50495  # was None

$ws.Range("H67").Value = 49999  # was 0
$ws.Range("J67").Value = 49999  # was 0
$ws.Range("L67").Value = 49999  # was 0
$ws.Range("N67").Value = -51715  # was None

$ws.Range("H102").Value = 485.35  # was 477.85715
$ws.Range("I102").Value = 485.35  # was 477.85715
$ws.Range("K102").Value = 485.35  # was 477.85715
$ws.Range("M102").Value = 1136.65  # was 1144.14285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 7710  # was 12800
$ws.Range("I3").Value = 2620  # was 0
$ws.Range("K3").Value = 2620  # was 0
$ws.Range("M3").Value = -2508  # was None

$ws.Range("H15").Value = 7710  # was 12800
$ws.Range("I15").Value = 2620  # was 0
$ws.Range("K15").Value = 2620  # was 0
$ws.Range("M15").Value = -2450  # was None

$ws.Range("H16").Value = 1860.8125  # was 1820.2941
$ws.Range("I16").Value = 1997.7693  # was 1995.4615
$ws.Range("J16").Value = 1267.3334  # was 1251
$ws.Range("K16").Value = 1997.7693  # was 1995.4615
$ws.Range("L16").Value = 1267.3334  # was 1251
$ws.Range("M16").Value = -1827.7693  # was -1825.4615
$ws.Range("N16").Value = -1607.3334  # was -1591

$ws.Range("H22").Value = 1500  # was 1495.8
$ws.Range("I22").Value = 0  # was 1487
$ws.Range("J22").Value = 1500  # was 1498
$ws.Range("K22").Value = 0  # was 1487
$ws.Range("L22").Value = 1500  # was 1498
$ws.Range("M22").ClearContents()  # was -1192
$ws.Range("N22").Value = -2090  # was -2088

$ws.Range("H27").Value = 1500  # was 1495.8
$ws.Range("I27").Value = 0  # was 1487
$ws.Range("J27").Value = 1500  # was 1498
$ws.Range("K27").Value = 0  # was 1487
$ws.Range("L27").Value = 1500  # was 1498
$ws.Range("M27").ClearContents()  # was -1380
$ws.Range("N27").Value = -1714  # was -1712

$ws.Range("H55").Value = 1529  # was 3079.8
$ws.Range("I55").Value = 1412.5  # was 1800
$ws.Range("J55").Value = 1595.5714  # was 4999.5
$ws.Range("K55").Value = 1412.5  # was 1800
$ws.Range("L55").Value = 1595.5714  # was 4999.5
$ws.Range("M55").Value = -1239.5  # was -1627
$ws.Range("N55").Value = -1941.5714  # was -5345.5

$ws.Range("H100").Value = 2000  # was 0
$ws.Range("I100").Value = 2000  # was 0
$ws.Range("K100").Value = 2000  # was 0
$ws.Range("M100").Value = -1459  # was None

$ws.Range("H122").Value = 3501.25  # was 3908.0908
$ws.Range("I122").Value = 3501.25  # was 3498.9
$ws.Range("J122").Value = 0  # was 8000
$ws.Range("K122").Value = 10503.75  # was 10496.7
$ws.Range("L122").Value = 0  # was 24000
$ws.Range("M122").Value = -8053.75  # was -8046.700000000001
$ws.Range("N122").ClearContents()  # was -28900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 0  # was 100000
$ws.Range("I138").Value = 0  # was 100000
$ws.Range("K138").Value = 0  # was 100000
$ws.Range("M138").ClearContents()  # was -94860
